# Update BOM value columns (A, C, E, G, I, K) per the pywrap ortools algorithm
# re-computation. Count columns (B, D, F, H, J, L) are left unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4
$ws.Range("C2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("G2").Value = 6
$ws.Range("I2").Value = 6
$ws.Range("K2").Value = 6
$ws.Range("A3").Value = 6
$ws.Range("C3").Value = 6
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("K3").Value = 8
$ws.Range("A4").Value = 6
$ws.Range("C4").Value = 6
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 12
$ws.Range("I4").Value = 12
$ws.Range("K4").Value = 12
$ws.Range("A5").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("E5").Value = 10
$ws.Range("G5").Value = 14
$ws.Range("I5").Value = 14
$ws.Range("K5").Value = 14
$ws.Range("A6").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("E6").Value = 10
$ws.Range("G6").Value = 16
$ws.Range("I6").Value = 16
$ws.Range("K6").Value = 16
$ws.Range("A7").Value = 11
$ws.Range("C7").Value = 11
$ws.Range("E7").Value = 11
$ws.Range("G7").Value = 20
$ws.Range("I7").Value = 20
$ws.Range("K7").Value = 20
$ws.Range("A8").Value = 13
$ws.Range("C8").Value = 13
$ws.Range("E8").Value = 13
$ws.Range("G8").Value = 22
$ws.Range("I8").Value = 22
$ws.Range("K8").Value = 22
$ws.Range("A9").Value = 15
$ws.Range("C9").Value = 15
$ws.Range("E9").Value = 15
$ws.Range("G9").Value = 41
$ws.Range("I9").Value = 41
$ws.Range("K9").Value = 41
$ws.Range("A10").Value = 32
$ws.Range("C10").Value = 32
$ws.Range("E10").Value = 32
$ws.Range("G10").Value = 48
$ws.Range("I10").Value = 48
$ws.Range("K10").Value = 48
$ws.Range("A11").Value = 41
$ws.Range("C11").Value = 41
$ws.Range("E11").Value = 41
$ws.Range("G11").Value = 60
$ws.Range("I11").Value = 60
$ws.Range("K11").Value = 60
$ws.Range("A12").Value = 42
$ws.Range("C12").Value = 42
$ws.Range("E12").Value = 42
$ws.Range("G12").Value = 62
$ws.Range("I12").Value = 62
$ws.Range("K12").Value = 62
$ws.Range("A13").Value = 46
$ws.Range("C13").Value = 46
$ws.Range("E13").Value = 46
$ws.Range("G13").Value = 80
$ws.Range("I13").Value = 80
$ws.Range("K13").Value = 80
$ws.Range("A14").Value = 51
$ws.Range("C14").Value = 51
$ws.Range("E14").Value = 51
$ws.Range("G14").Value = 86
$ws.Range("I14").Value = 86
$ws.Range("K14").Value = 86
$ws.Range("A15").Value = 62
$ws.Range("C15").Value = 62
$ws.Range("E15").Value = 62
$ws.Range("G15").Value = 93
$ws.Range("I15").Value = 93
$ws.Range("K15").Value = 93
$ws.Range("A16").Value = 82
$ws.Range("C16").Value = 82
$ws.Range("E16").Value = 82
$ws.Range("G16").Value = 96
$ws.Range("I16").Value = 96
$ws.Range("K16").Value = 96
$ws.Range("A17").Value = 86
$ws.Range("C17").Value = 86
$ws.Range("E17").Value = 86
$ws.Range("G17").Value = 120
$ws.Range("I17").Value = 120
$ws.Range("A18").Value = 93
$ws.Range("C18").Value = 93
$ws.Range("E18").Value = 93
$ws.Range("A19").Value = 96
$ws.Range("C19").Value = 96
$ws.Range("E19").Value = 96
